# Auto-generated Excel COM-interop script
# Applies updated market-price values (columns H-N) across multiple sheets
# as produced by the scheduled Sheets runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1101.8572
$ws.Range("I6").Value = 838.6667
$ws.Range("K6").Value = 2516.0001
$ws.Range("M6").Value = -2404.0001

$ws.Range("H18").Value = 5185.1665
$ws.Range("I18").Value = 5371.778
$ws.Range("K18").Value = 5371.778
$ws.Range("M18").Value = -5087.778

$ws.Range("H43").Value = 2480.3333
$ws.Range("J43").Value = 2480
$ws.Range("L43").Value = 2480
$ws.Range("N43").Value = -2618

$ws.Range("H53").Value = 1984.5
$ws.Range("J53").Value = 1966.3334
$ws.Range("L53").Value = 1966.3334
$ws.Range("N53").Value = -3240.3334

$ws.Range("H115").Value = 623.6
$ws.Range("I115").Value = 623.6
$ws.Range("K115").Value = 1870.8
$ws.Range("M115").Value = -303.8000000000002

$ws.Range("H127").Value = 4475.5
$ws.Range("I127").Value = 4467
$ws.Range("J127").Value = 4492.5
$ws.Range("K127").Value = 13401
$ws.Range("L127").Value = 13477.5
$ws.Range("M127").Value = -8441
$ws.Range("N127").Value = -23397.5

$ws.Range("H131").Value = 8822.259
$ws.Range("I131").Value = 7777.227
$ws.Range("K131").Value = 23331.681
$ws.Range("M131").Value = -18291.681

$ws.Range("H132").Value = 3460.4783
$ws.Range("I132").Value = 3533.682
$ws.Range("K132").Value = 10601.046
$ws.Range("M132").Value = -8071.045999999998

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").ClearContents()
$ws.Range("N134").Value = 0

$ws.Range("H135").Value = 2989
$ws.Range("I135").Value = 2989
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 26901
$ws.Range("L135").Value = 0
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -24366

$ws.Range("H141").Value = 4336.5713
$ws.Range("I141").Value = 3148.2104
$ws.Range("J141").Value = 15626
$ws.Range("K141").Value = 9444.6312
$ws.Range("L141").Value = 46878
$ws.Range("M141").Value = -4264.6312
$ws.Range("N141").Value = -57238

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4421.2544
$ws.Range("I32").Value = 4536.325
$ws.Range("K32").Value = 4536.325
$ws.Range("M32").Value = -4249.325

$ws.Range("H44").Value = 39043.5
$ws.Range("J44").Value = 39043.5
$ws.Range("L44").Value = 39043.5
$ws.Range("N44").Value = -40019.5

$ws.Range("H80").Value = 54992.5
$ws.Range("J80").Value = 54992.5
$ws.Range("L80").Value = 54992.5
$ws.Range("N80").Value = -56988.5

$ws.Range("H83").Value = 54992.5
$ws.Range("J83").Value = 54992.5
$ws.Range("L83").Value = 164977.5
$ws.Range("N83").Value = -174961.5

$ws.Range("H132").Value = 2582.6943
$ws.Range("I132").Value = 2579.6
$ws.Range("K132").Value = 7738.799999999999
$ws.Range("M132").Value = -5208.799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1769.9231
$ws.Range("I99").Value = 1872.5454
$ws.Range("K99").Value = 1872.5454
$ws.Range("M99").Value = -374.5454

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 34461.3
$ws.Range("J41").Value = 49999
$ws.Range("L41").Value = 49999
$ws.Range("N41").Value = -50855

$ws.Range("H50").Value = 43329.668
$ws.Range("J50").Value = 49994.5
$ws.Range("L50").Value = 49994.5
$ws.Range("N50").Value = -51244.5

$ws.Range("H51").Value = 47000
$ws.Range("J51").Value = 47000
$ws.Range("L51").Value = 47000
$ws.Range("N51").Value = -48472

$ws.Range("H60").Value = 22765.6
$ws.Range("J60").Value = 32075.4
$ws.Range("L60").Value = 32075.4
$ws.Range("N60").Value = -33097.4

$ws.Range("H61").Value = 47000
$ws.Range("J61").Value = 47000
$ws.Range("L61").Value = 47000
$ws.Range("N61").Value = -47696

$ws.Range("H68").Value = 43692.6
$ws.Range("J68").Value = 67981.5
$ws.Range("L68").Value = 67981.5
$ws.Range("N68").Value = -69479.5

$ws.Range("H71").Value = 43692.6
$ws.Range("J71").Value = 67981.5
$ws.Range("L71").Value = 203944.5
$ws.Range("N71").Value = -211432.5

$ws.Range("H74").Value = 69916.164
$ws.Range("J74").Value = 69916.164
$ws.Range("L74").Value = 69916.164
$ws.Range("N74").Value = -71664.164

$ws.Range("H77").Value = 69916.164
$ws.Range("J77").Value = 69916.164
$ws.Range("L77").Value = 209748.492
$ws.Range("N77").Value = -218484.492

$ws.Range("H107").Value = 35748476
$ws.Range("I107").Value = 45496260
$ws.Range("J107").Value = 6600.1665
$ws.Range("K107").Value = 45496260
$ws.Range("L107").Value = 6600.1665
$ws.Range("M107").Value = -45494340
$ws.Range("N107").Value = -10440.1665

$ws.Range("H122").Value = 1765.3334
$ws.Range("I122").Value = 1318.9286
$ws.Range("J122").Value = 2390.3
$ws.Range("K122").Value = 3956.7858
$ws.Range("L122").Value = 7170.900000000001
$ws.Range("M122").Value = -1506.7858
$ws.Range("N122").Value = -12070.9

$ws.Range("H132").Value = 7213.2573
$ws.Range("I132").Value = 2713.0952
$ws.Range("J132").Value = 13963.5
$ws.Range("K132").Value = 8139.285600000001
$ws.Range("L132").Value = 41890.5
$ws.Range("M132").Value = -5609.285600000001
$ws.Range("N132").Value = -46950.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()

$ws.Range("H103").Value = 8974.444
$ws.Range("I103").Value = 9352.857
$ws.Range("J103").Value = 7650
$ws.Range("K103").Value = 28058.571
$ws.Range("L103").Value = 22950
$ws.Range("M103").Value = -27179.571
$ws.Range("N103").Value = -24708

$ws.Range("H124").Value = 5003.3335
$ws.Range("J124").Value = 5578.857
$ws.Range("L124").Value = 16736.571
$ws.Range("N124").Value = -26556.571

$ws.Range("H137").Value = 2364.5874
$ws.Range("J137").Value = 2448.8904
$ws.Range("L137").Value = 7346.671200000001
$ws.Range("N137").Value = -17546.6712

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3374
$ws.Range("I80").Value = 3166.3333
$ws.Range("K80").Value = 3166.3333
$ws.Range("M80").Value = -2168.3333

$ws.Range("H83").Value = 3374
$ws.Range("I83").Value = 3166.3333
$ws.Range("K83").Value = 15831.6665
$ws.Range("M83").Value = -10839.6665

$ws.Range("H93").Value = 59500
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()

$ws.Range("H132").Value = 3881.3333
$ws.Range("I132").Value = 3881.3333
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11643.9999
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -9113.999899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 1441.6
$ws.Range("I13").Value = 686.3333
$ws.Range("J13").Value = 2574.5
$ws.Range("K13").Value = 686.3333
$ws.Range("L13").Value = 2574.5
$ws.Range("M13").Value = -546.3333
$ws.Range("N13").Value = -2854.5

$ws.Range("H46").Value = 5396.5713
$ws.Range("J46").Value = 5396.5713
$ws.Range("L46").Value = 5396.5713
$ws.Range("N46").Value = -5772.5713

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 24325.334
$ws.Range("I70").Value = 22988
$ws.Range("K70").Value = 22988
$ws.Range("M70").Value = -22673

$ws.Range("H73").Value = 24325.334
$ws.Range("I73").Value = 22988
$ws.Range("K73").Value = 22988
$ws.Range("M73").Value = -21896
